$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C3"
$ws.Range("C2").Value = "Cd19"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.857161
$ws.Range("H2").Value = 5.571483000000001
$ws.Range("I2").Value = 0.006673232049902625
$ws.Range("J2").Value = 0.006673232049902625
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.284271
$ws.Range("N2").Value = 0.852813
$ws.Range("O2").Value = 0.9612929915054016
$ws.Range("P2").Value = 0.9612929915054016
$ws.Range("Q2").Value = 0.5279370146310001
$ws.Range("R2").Value = 4.751433131679001
$ws.Range("S2").Value = 0.006414931200260618
$ws.Range("T2").Value = 0.006414931200260618

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C3"
$ws.Range("C3").Value = "Cd19"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.857161
$ws.Range("H3").Value = 5.571483000000001
$ws.Range("I3").Value = 0.006673232049902625
$ws.Range("J3").Value = 0.006673232049902625
$ws.Range("K3").Value = 1.0
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01144633333333333
$ws.Range("N3").Value = 0.034339
$ws.Range("O3").Value = 0.03870700849459845
$ws.Range("P3").Value = 0.03870700849459845
$ws.Range("Q3").Value = 0.02125768385966667
$ws.Range("R3").Value = 0.191319154737
$ws.Range("S3").Value = 0.0002583008496420075
$ws.Range("T3").Value = 0.0002583008496420075

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "C3"
$ws.Range("C4").Value = "Cd19"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 82.87880433333333
$ws.Range("H4").Value = 248.636413
$ws.Range("I4").Value = 0.2978037409437354
$ws.Range("J4").Value = 0.2978037409437354
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.284271
$ws.Range("N4").Value = 0.852813
$ws.Range("O4").Value = 0.9612929915054016
$ws.Range("P4").Value = 0.9612929915054016
$ws.Range("Q4").Value = 23.560040586641
$ws.Range("R4").Value = 212.040365279769
$ws.Range("S4").Value = 0.286276649013303
$ws.Range("T4").Value = 0.286276649013303

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C3"
$ws.Range("C5").Value = "Cd19"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 82.87880433333333
$ws.Range("H5").Value = 248.636413
$ws.Range("I5").Value = 0.2978037409437354
$ws.Range("J5").Value = 0.2978037409437354
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01144633333333333
$ws.Range("N5").Value = 0.034339
$ws.Range("O5").Value = 0.03870700849459845
$ws.Range("P5").Value = 0.03870700849459845
$ws.Range("Q5").Value = 0.9486584206674445
$ws.Range("R5").Value = 8.537925786007001
$ws.Range("S5").Value = 0.01152709193043236
$ws.Range("T5").Value = 0.01152709193043236

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "C3"
$ws.Range("C6").Value = "Cd19"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 59.75754533333333
$ws.Range("H6").Value = 179.272636
$ws.Range("I6").Value = 0.214723422870666
$ws.Range("J6").Value = 0.214723422870666
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.284271
$ws.Range("N6").Value = 0.852813
$ws.Range("O6").Value = 0.9612929915054016
$ws.Range("P6").Value = 0.9612929915054016
$ws.Range("Q6").Value = 16.987337169452
$ws.Range("R6").Value = 152.886034525068
$ws.Range("S6").Value = 0.2064121215176219
$ws.Range("T6").Value = 0.2064121215176219

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "C3"
$ws.Range("C7").Value = "Cd19"
$ws.Range("D7").Value = "Neutrophils"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 59.75754533333333
$ws.Range("H7").Value = 179.272636
$ws.Range("I7").Value = 0.214723422870666
$ws.Range("J7").Value = 0.214723422870666
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01144633333333333
$ws.Range("N7").Value = 0.034339
$ws.Range("O7").Value = 0.03870700849459845
$ws.Range("P7").Value = 0.03870700849459845
$ws.Range("Q7").Value = 0.6840047830671112
$ws.Range("R7").Value = 6.156043047604
$ws.Range("S7").Value = 0.008311301353044124
$ws.Range("T7").Value = 0.008311301353044125

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "C3"
$ws.Range("C8").Value = "Cd19"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 0.3694876666666667
$ws.Range("H8").Value = 1.108463
$ws.Range("I8").Value = 0.00132765922784494
$ws.Range("J8").Value = 0.00132765922784494
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.284271
$ws.Range("N8").Value = 0.852813
$ws.Range("O8").Value = 0.9612929915054016
$ws.Range("P8").Value = 0.9612929915054016
$ws.Range("Q8").Value = 0.105034628491
$ws.Range("R8").Value = 0.9453116564190001
$ws.Range("S8").Value = 0.001276269510834814
$ws.Range("T8").Value = 0.001276269510834814

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "C3"
$ws.Range("C9").Value = "Cd19"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 0.3694876666666667
$ws.Range("H9").Value = 1.108463
$ws.Range("I9").Value = 0.00132765922784494
$ws.Range("J9").Value = 0.00132765922784494
$ws.Range("K9").Value = 1.0
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01144633333333333
$ws.Range("N9").Value = 0.034339
$ws.Range("O9").Value = 0.03870700849459845
$ws.Range("P9").Value = 0.03870700849459845
$ws.Range("Q9").Value = 0.004229278995222223
$ws.Range("R9").Value = 0.038063510957
$ws.Range("S9").Value = 0.00005138971701012612
$ws.Range("T9").Value = 0.00005138971701012612

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "C3"
$ws.Range("C10").Value = "Cd19"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 127.8666333333333
$ws.Range("H10").Value = 383.5999
$ws.Range("I10").Value = 0.4594559737541049
$ws.Range("J10").Value = 0.4594559737541049
$ws.Range("K10").Value = 1.0
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.284271
$ws.Range("N10").Value = 0.852813
$ws.Range("O10").Value = 0.9612929915054016
$ws.Range("P10").Value = 0.9612929915054016
$ws.Range("Q10").Value = 36.3487757243
$ws.Range("R10").Value = 327.1389815187
$ws.Range("S10").Value = 0.4416718074751108
$ws.Range("T10").Value = 0.4416718074751108

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "C3"
$ws.Range("C11").Value = "Cd19"
$ws.Range("D11").Value = "Neutrophils"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 127.8666333333333
$ws.Range("H11").Value = 383.5999
$ws.Range("I11").Value = 0.4594559737541049
$ws.Range("J11").Value = 0.4594559737541049
$ws.Range("K11").Value = 1.0
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01144633333333333
$ws.Range("N11").Value = 0.034339
$ws.Range("O11").Value = 0.03870700849459845
$ws.Range("P11").Value = 0.03870700849459845
$ws.Range("Q11").Value = 1.463604107344445
$ws.Range("R11").Value = 13.1724369661
$ws.Range("S11").Value = 0.01778416627899414
$ws.Range("T11").Value = 0.01778416627899414

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "C3"
$ws.Range("C12").Value = "Cd19"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 5.570446333333334
$ws.Range("H12").Value = 16.711339
$ws.Range("I12").Value = 0.02001597115374626
$ws.Range("J12").Value = 0.02001597115374626
$ws.Range("K12").Value = 1.0
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.284271
$ws.Range("N12").Value = 0.852813
$ws.Range("O12").Value = 0.9612929915054016
$ws.Range("P12").Value = 0.9612929915054016
$ws.Range("Q12").Value = 1.583516349623
$ws.Range("R12").Value = 14.251647146607
$ws.Range("S12").Value = 0.01924121278827057
$ws.Range("T12").Value = 0.01924121278827057

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "C3"
$ws.Range("C13").Value = "Cd19"
$ws.Range("D13").Value = "Neutrophils"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 5.570446333333334
$ws.Range("H13").Value = 16.711339
$ws.Range("I13").Value = 0.02001597115374626
$ws.Range("J13").Value = 0.02001597115374626
$ws.Range("K13").Value = 1.0
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.01144633333333333
$ws.Range("N13").Value = 0.034339
$ws.Range("O13").Value = 0.03870700849459845
$ws.Range("P13").Value = 0.03870700849459845
$ws.Range("Q13").Value = 0.06376118554677779
$ws.Range("R13").Value = 0.5738506699210001
$ws.Range("S13").Value = 0.0007747583654756939
$ws.Range("T13").Value = 0.000774758365475694
